$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A186").Value = 185
$ws.Range("B186").Value = 1
$ws.Range("C186").Value = "2024-06-18 21:11:47"
$ws.Range("D186").Value = 200
$ws.Range("E186").Value = 13

$ws.Range("A187").Value = 186
$ws.Range("B187").Value = 2
$ws.Range("C187").Value = "2024-06-18 21:11:48"
$ws.Range("D187").Value = 200
$ws.Range("E187").Value = 3
